# Refresh the crypto price/volume snapshot (GitHub Actions bot run).
# Column D (Price) and column E (Volume(1h)) are stored as plain text in
# this sheet (e.g. "3.440.43", "  +0.54%  "), so every write below keeps the
# literal text instead of letting Excel auto-coerce numeric-looking prices
# into actual numbers (which would silently drop things like trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.762.89"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "3.440.43"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value2 = "'575.80"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value2 = "'160.28"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.441.36"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value2 = "'0.582"
$ws.Range("E9").Value = "  +8.88%  "
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value2 = "'0.441"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "4.035.62"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value2 = "'28.28"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").Value = "64.771.41"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "3.429.85"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value2 = "'6.37"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value2 = "'14.29"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value2 = "'387.01"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value2 = "'8.18"
$ws.Range("E22").Value = "  -3.63%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  +14.94%  "
$ws.Range("D27").Value2 = "'9.81"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").Value2 = "'0.180"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value2 = "'0.999"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value2 = "'6.19"
$ws.Range("E30").Value = "  +6.87%  "
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value2 = "'6.57"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").Value2 = "'23.68"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value2 = "'7.10"
$ws.Range("E36").Value = "  +3.72%  "
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").Value2 = "'163.34"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").Value = "3.020.58"
$ws.Range("E39").Value = "  +5.22%  "
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").Value2 = "'27.22"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("D43").Value2 = "'4.55"
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("D44").Value2 = "'42.79"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").Value2 = "'0.0317"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").Value2 = "'24.68"
$ws.Range("D48").Value2 = "'1.09"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value2 = "'0.878"
$ws.Range("E49").Value = "  +6.16%  "
$ws.Range("D50").Value2 = "'6.64"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").Value2 = "'2.18"
$ws.Range("E51").Value = "  +4.40%  "
